$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that is guaranteed to land in the sheet as TEXT
# (shared string), even when it looks like a number (e.g. "0.675080691881").
# A direct `.Value = "0.675..."` assignment gets auto-coerced to a numeric
# cell by Excel, so instead we briefly place a text-producing formula in the
# cell, then convert it in place to a static value via copy / paste-values.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = 0
}

# --- Header row (row 1) ---
$ws.Range("A1").Value = "METHOD"
$ws.Range("B1").Value = "HIGHEST ACCURACY"
$ws.Range("C1").Value = "NORM?"
$ws.Range("D1").Value = "PARAM (non-default)"

# --- Row 2: KNN ---
$ws.Range("A2").Value = "KNN"
Set-TextValue $ws.Range("B2") " 0.648183398163"
$ws.Range("C2").Value = "normalized"
$ws.Range("D2").Value = "n_neighbors=10, weights= 'distance',metric='manhattan'"

# --- Row 3: RANDOM FOREST ---
$ws.Range("A3").Value = "RANDOM FOREST"
Set-TextValue $ws.Range("B3") "0.67218406025"
$ws.Range("D3").Value = 'n_estimators=250, max_depth=None, bootstrap=False, class_weight="balanced", n_jobs=4'

# --- Row 4: DECISION TREE ---
$ws.Range("A4").Value = "DECISION TREE"
Set-TextValue $ws.Range("B4") " 0.5796"
$ws.Range("D4").Value = "max_depth=10, splitter='best', min_samples_split=81"

# --- Row 5: LOGISTIC REG ---
$ws.Range("A5").Value = "LOGISTIC REG"
$ws.Range("B5").Value = "?????????"
$ws.Range("D5").Value = 'multi_class=''multinomial'',  max_iter=500, solver="newton-cg",C=1'

# --- Row 6: SVM ---
$ws.Range("A6").Value = "SVM"
Set-TextValue $ws.Range("B6") "0.545973682033"
$ws.Range("C6").Value = "normalized"
$ws.Range("D6").Value = "C=1000, kernel='poly', degree=10, coef0=3"

# --- Row 7: NEURAL NETWORKS ---
$ws.Range("A7").Value = "NEURAL NETWORKS"

# --- Row 8: Ensemble ---
$ws.Range("A8").Value = "Ensemble"
Set-TextValue $ws.Range("B8") "0.675080691881"

# --- Apply "Bad" cell style (red text / pink fill) to the two flagged rows ---
$ws.Range("A5").Style = "Bad"
$ws.Range("A7").Style = "Bad"

# --- Selection change to match the saved view state ---
$ws.Range("B12").Select()
